$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meandam (E), lowerdam (F), upperdam (G) damage values per row
$damageValues = @{
    2 = @(786974946.6503377, 590231209.9877533, 983718683.3129222)
    3 = @(8556928.978965666, 6417696.73422425, 10696161.22370708)
    4 = @(822769826.0269077, 617077369.5201807, 1028462282.533635)
    5 = @(8840563.944195945, 6630422.958146958, 11050704.93024493)
    6 = @(906710923.0051646, 680033192.2538733, 1133388653.756456)
    7 = @(9361319.396079145, 7020989.547059359, 11701649.24509893)
    8 = @(961991683.7428241, 721493762.8071179, 1202489604.67853)
    9 = @(9607635.674413536, 7205726.755810153, 12009544.59301692)
    10 = @(1048718331.973319, 786538748.9799889, 1310897914.966648)
    11 = @(9907322.48875894, 7430491.866569204, 12384153.11094867)
    12 = @(1108623245.814198, 831467434.3606486, 1385779057.267748)
    13 = @(10152154.5233175, 7614115.892488128, 12690193.15414688)
    14 = @(1165745662.727023, 874309247.0452669, 1457182078.408778)
    15 = @(10412356.72956229, 7809267.547171716, 13015445.91195286)
    16 = @(1239797163.893575, 929847872.9201813, 1549746454.866969)
    17 = @(10917438.68401035, 8188079.01300776, 13646798.35501293)
    18 = @(1299836983.490195, 974877737.6176465, 1624796229.362744)
    19 = @(11390180.25532855, 8542635.191496413, 14237725.31916069)
    20 = @(1362325886.008341, 1021744414.506256, 1702907357.510426)
    21 = @(11851881.3256731, 8888910.99425482, 14814851.65709137)
}

foreach ($row in $damageValues.Keys) {
    $vals = $damageValues[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]  # column E - meandam
    $ws.Cells.Item($row, 6).Value = $vals[1]  # column F - lowerdam
    $ws.Cells.Item($row, 7).Value = $vals[2]  # column G - upperdam
}
